$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows at position 166, shifting existing rows 166-191 down to 170-195
$ws.Rows("166:169").Insert()

# Temporarily mark the date (B) and id (C) columns as text so values such as
# "2019-11-18" and "03017" are written as literal text rather than being auto-converted
# to a date serial number / plain integer by Excel's type inference.
$ws.Range("B166:C195").NumberFormat = "@"

# Write final values for rows 166-195 (new rows + re-ordered/updated existing rows)
$ws.Range("A166").Value = 1574035200
$ws.Range("B166").Value = "2019-11-18"
$ws.Range("C166").Value = "03017"
$ws.Range("D166").Value = "UNIWALL"
$ws.Range("E166").Value = 0.98
$ws.Range("F166").Value = 0.98
$ws.Range("G166").Value = 0.98
$ws.Range("H166").Value = 0.98
$ws.Range("I166").Value = 10000

$ws.Range("A167").Value = 1574380800
$ws.Range("B167").Value = "2019-11-22"
$ws.Range("C167").Value = "03017"
$ws.Range("D167").Value = "UNIWALL"
$ws.Range("E167").Value = 1
$ws.Range("F167").Value = 1
$ws.Range("G167").Value = 1
$ws.Range("H167").Value = 1
$ws.Range("I167").Value = 10000

$ws.Range("A168").Value = 1574812800
$ws.Range("B168").Value = "2019-11-27"
$ws.Range("C168").Value = "03017"
$ws.Range("D168").Value = "UNIWALL"
$ws.Range("E168").Value = 1.02
$ws.Range("F168").Value = 1.02
$ws.Range("G168").Value = 1.02
$ws.Range("H168").Value = 1.02
$ws.Range("I168").Value = 20000

$ws.Range("A169").Value = 1574899200
$ws.Range("B169").Value = "2019-11-28"
$ws.Range("C169").Value = "03017"
$ws.Range("D169").Value = "UNIWALL"
$ws.Range("E169").Value = 1.04
$ws.Range("F169").Value = 1.05
$ws.Range("G169").Value = 1.04
$ws.Range("H169").Value = 1.05
$ws.Range("I169").Value = 30000

$ws.Range("A170").Value = 1574985600
$ws.Range("B170").Value = "2019-11-29"
$ws.Range("C170").Value = "03017"
$ws.Range("D170").Value = "UNIWALL"
$ws.Range("E170").Value = 1.05
$ws.Range("F170").Value = 1.05
$ws.Range("G170").Value = 1.05
$ws.Range("H170").Value = 1.05
$ws.Range("I170").Value = "-"

$ws.Range("A171").Value = 1575417600
$ws.Range("B171").Value = "2019-12-04"
$ws.Range("C171").Value = "03017"
$ws.Range("D171").Value = "UNIWALL"
$ws.Range("E171").Value = 1.05
$ws.Range("F171").Value = 1.05
$ws.Range("G171").Value = 1.05
$ws.Range("H171").Value = 1.05
$ws.Range("I171").Value = "-"

$ws.Range("A172").Value = 1575849600
$ws.Range("B172").Value = "2019-12-09"
$ws.Range("C172").Value = "03017"
$ws.Range("D172").Value = "UNIWALL"
$ws.Range("E172").Value = 1.06
$ws.Range("F172").Value = 1.09
$ws.Range("G172").Value = 1.06
$ws.Range("H172").Value = 1.09
$ws.Range("I172").Value = 70000

$ws.Range("A173").Value = 1575936000
$ws.Range("B173").Value = "2019-12-10"
$ws.Range("C173").Value = "03017"
$ws.Range("D173").Value = "UNIWALL"
$ws.Range("E173").Value = 1.1
$ws.Range("F173").Value = 1.1
$ws.Range("G173").Value = 1.1
$ws.Range("H173").Value = 1.1
$ws.Range("I173").Value = 20000

$ws.Range("A174").Value = 1576022400
$ws.Range("B174").Value = "2019-12-11"
$ws.Range("C174").Value = "03017"
$ws.Range("D174").Value = "UNIWALL"
$ws.Range("E174").Value = 1.12
$ws.Range("F174").Value = 1.13
$ws.Range("G174").Value = 1.12
$ws.Range("H174").Value = 1.13
$ws.Range("I174").Value = 13000

$ws.Range("A175").Value = 1576108800
$ws.Range("B175").Value = "2019-12-12"
$ws.Range("C175").Value = "03017"
$ws.Range("D175").Value = "UNIWALL"
$ws.Range("E175").Value = 1.14
$ws.Range("F175").Value = 1.14
$ws.Range("G175").Value = 1.14
$ws.Range("H175").Value = 1.14
$ws.Range("I175").Value = 10000

$ws.Range("A176").Value = 1576195200
$ws.Range("B176").Value = "2019-12-13"
$ws.Range("C176").Value = "03017"
$ws.Range("D176").Value = "UNIWALL"
$ws.Range("E176").Value = 1.15
$ws.Range("F176").Value = 1.15
$ws.Range("G176").Value = 1.15
$ws.Range("H176").Value = 1.15
$ws.Range("I176").Value = 20000

$ws.Range("A177").Value = 1576627200
$ws.Range("B177").Value = "2019-12-18"
$ws.Range("C177").Value = "03017"
$ws.Range("D177").Value = "UNIWALL"
$ws.Range("E177").Value = 1.15
$ws.Range("F177").Value = 1.15
$ws.Range("G177").Value = 1.15
$ws.Range("H177").Value = 1.15
$ws.Range("I177").Value = "-"

$ws.Range("A178").Value = 1576713600
$ws.Range("B178").Value = "2019-12-19"
$ws.Range("C178").Value = "03017"
$ws.Range("D178").Value = "UNIWALL"
$ws.Range("E178").Value = 1.15
$ws.Range("F178").Value = 1.15
$ws.Range("G178").Value = 1.15
$ws.Range("H178").Value = 1.15
$ws.Range("I178").Value = "-"

$ws.Range("A179").Value = 1576800000
$ws.Range("B179").Value = "2019-12-20"
$ws.Range("C179").Value = "03017"
$ws.Range("D179").Value = "UNIWALL"
$ws.Range("E179").Value = 1.15
$ws.Range("F179").Value = 1.15
$ws.Range("G179").Value = 1.15
$ws.Range("H179").Value = 1.15
$ws.Range("I179").Value = "-"

$ws.Range("A180").Value = 1577059200
$ws.Range("B180").Value = "2019-12-23"
$ws.Range("C180").Value = "03017"
$ws.Range("D180").Value = "UNIWALL"
$ws.Range("E180").Value = 1.18
$ws.Range("F180").Value = 1.2
$ws.Range("G180").Value = 1.18
$ws.Range("H180").Value = 1.2
$ws.Range("I180").Value = 20000

$ws.Range("A181").Value = 1578873600
$ws.Range("B181").Value = "2020-01-13"
$ws.Range("C181").Value = "03017"
$ws.Range("D181").Value = "UNIWALL"
$ws.Range("E181").Value = 1.2
$ws.Range("F181").Value = 1.2
$ws.Range("G181").Value = 1.2
$ws.Range("H181").Value = 1.2
$ws.Range("I181").Value = 200000

$ws.Range("A182").Value = 1579046400
$ws.Range("B182").Value = "2020-01-15"
$ws.Range("C182").Value = "03017"
$ws.Range("D182").Value = "UNIWALL"
$ws.Range("E182").Value = 1.2
$ws.Range("F182").Value = 1.2
$ws.Range("G182").Value = 1.2
$ws.Range("H182").Value = 1.2
$ws.Range("I182").Value = "-"

$ws.Range("A183").Value = 1579132800
$ws.Range("B183").Value = "2020-01-16"
$ws.Range("C183").Value = "03017"
$ws.Range("D183").Value = "UNIWALL"
$ws.Range("E183").Value = 1.2
$ws.Range("F183").Value = 1.2
$ws.Range("G183").Value = 1.2
$ws.Range("H183").Value = 1.2
$ws.Range("I183").Value = "-"

$ws.Range("A184").Value = 1580256000
$ws.Range("B184").Value = "2020-01-29"
$ws.Range("C184").Value = "03017"
$ws.Range("D184").Value = "UNIWALL"
$ws.Range("E184").Value = 1.21
$ws.Range("F184").Value = 1.21
$ws.Range("G184").Value = 1.21
$ws.Range("H184").Value = 1.21
$ws.Range("I184").Value = 30000

$ws.Range("A185").Value = 1580342400
$ws.Range("B185").Value = "2020-01-30"
$ws.Range("C185").Value = "03017"
$ws.Range("D185").Value = "UNIWALL"
$ws.Range("E185").Value = 1.21
$ws.Range("F185").Value = 1.21
$ws.Range("G185").Value = 1.21
$ws.Range("H185").Value = 1.21
$ws.Range("I185").Value = "-"

$ws.Range("A186").Value = 1580428800
$ws.Range("B186").Value = "2020-01-31"
$ws.Range("C186").Value = "03017"
$ws.Range("D186").Value = "UNIWALL"
$ws.Range("E186").Value = 1.21
$ws.Range("F186").Value = 1.21
$ws.Range("G186").Value = 1.21
$ws.Range("H186").Value = 1.21
$ws.Range("I186").Value = "-"

$ws.Range("A187").Value = 1581638400
$ws.Range("B187").Value = "2020-02-14"
$ws.Range("C187").Value = "03017"
$ws.Range("D187").Value = "UNIWALL"
$ws.Range("E187").Value = 1.22
$ws.Range("F187").Value = 1.22
$ws.Range("G187").Value = 1.22
$ws.Range("H187").Value = 1.22
$ws.Range("I187").Value = 12000

$ws.Range("A188").Value = 1583107200
$ws.Range("B188").Value = "2020-03-02"
$ws.Range("C188").Value = "03017"
$ws.Range("D188").Value = "UNIWALL"
$ws.Range("E188").Value = 1.22
$ws.Range("F188").Value = 1.22
$ws.Range("G188").Value = 1.22
$ws.Range("H188").Value = 1.22
$ws.Range("I188").Value = "-"

$ws.Range("A189").Value = 1583193600
$ws.Range("B189").Value = "2020-03-03"
$ws.Range("C189").Value = "03017"
$ws.Range("D189").Value = "UNIWALL"
$ws.Range("E189").Value = 1.3
$ws.Range("F189").Value = 1.3
$ws.Range("G189").Value = 1.3
$ws.Range("H189").Value = 1.3
$ws.Range("I189").Value = 5000

$ws.Range("A190").Value = 1583280000
$ws.Range("B190").Value = "2020-03-04"
$ws.Range("C190").Value = "03017"
$ws.Range("D190").Value = "UNIWALL"
$ws.Range("E190").Value = 1.35
$ws.Range("F190").Value = 1.35
$ws.Range("G190").Value = 1.35
$ws.Range("H190").Value = 1.35
$ws.Range("I190").Value = 10000

$ws.Range("A191").Value = 1583712000
$ws.Range("B191").Value = "2020-03-09"
$ws.Range("C191").Value = "03017"
$ws.Range("D191").Value = "UNIWALL"
$ws.Range("E191").Value = 1.35
$ws.Range("F191").Value = 1.35
$ws.Range("G191").Value = 1.35
$ws.Range("H191").Value = 1.35
$ws.Range("I191").Value = "-"

$ws.Range("A192").Value = 1583798400
$ws.Range("B192").Value = "2020-03-10"
$ws.Range("C192").Value = "03017"
$ws.Range("D192").Value = "UNIWALL"
$ws.Range("E192").Value = 1.35
$ws.Range("F192").Value = 1.35
$ws.Range("G192").Value = 1.35
$ws.Range("H192").Value = 1.35
$ws.Range("I192").Value = "-"

$ws.Range("A193").Value = 1583884800
$ws.Range("B193").Value = "2020-03-11"
$ws.Range("C193").Value = "03017"
$ws.Range("D193").Value = "UNIWALL"
$ws.Range("E193").Value = 1.35
$ws.Range("F193").Value = 1.35
$ws.Range("G193").Value = 1.35
$ws.Range("H193").Value = 1.35
$ws.Range("I193").Value = "-"

$ws.Range("A194").Value = 1583971200
$ws.Range("B194").Value = "2020-03-12"
$ws.Range("C194").Value = "03017"
$ws.Range("D194").Value = "UNIWALL"
$ws.Range("E194").Value = 1.35
$ws.Range("F194").Value = 1.35
$ws.Range("G194").Value = 1.35
$ws.Range("H194").Value = 1.35
$ws.Range("I194").Value = "-"

$ws.Range("A195").Value = 1584057600
$ws.Range("B195").Value = "2020-03-13"
$ws.Range("C195").Value = "03017"
$ws.Range("D195").Value = "UNIWALL"
$ws.Range("E195").Value = 1.35
$ws.Range("F195").Value = 1.35
$ws.Range("G195").Value = 1.35
$ws.Range("H195").Value = 1.35
$ws.Range("I195").Value = "-"

# Remove the temporary text formatting now that the literal values are stored,
# restoring the cells to the same unstyled state as the rest of the data rows.
$ws.Range("B166:C195").ClearFormats()